$wb = $excel.ActiveWorkbook

# --- PageObjectModel: update/extend the Action key words list ---
$wsPage = $wb.Worksheets.Item("PageObjectModel")
$wsPage.Range("B22").Value = "waitTextToBePresented"
$wsPage.Range("B23").Value = "waitForCheck"
$wsPage.Range("B24").Value = "waitUtilSelected"
$wsPage.Range("B25").Value = "clickAndHold"
$wsPage.Range("B26").Value = "dragAndDropByOffset"
$wsPage.Range("B27").Value = "selectPartialContent"

# --- RichTextTestPage: fix the Target Name value ---
$wsRich = $wb.Worksheets.Item("RichTextTestPage")
$wsRich.Range("C2").Value = "RichTextBox"

# --- Selection bookkeeping on PageObjectModel (view state only) ---
$wsPage.Range("E19").Select()

# --- Switch the active sheet to RichTextTestPage and select row 3 ---
$wsRich.Activate()
$wsRich.Range("A3:XFD3").Select()
